$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.395.89'
$ws.Range("E2").Value = '  +3.29%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.789.43'
$ws.Range("E3").Value = '  +3.52%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.003'
$ws.Range("E4").Value = '  +0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '335.90'
$ws.Range("E5").Value = '  +0.87%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  +0.07%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3788'
$ws.Range("E7").Value = '  +1.22%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3423'
$ws.Range("E8").Value = '  +1.20%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '48.00'
$ws.Range("E9").Value = '  -1.12%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.202'
$ws.Range("E10").Value = '  +1.48%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07471'
$ws.Range("E11").Value = '  -0.10%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.001'
$ws.Range("E12").Value = '  -0.05%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '21.98'
$ws.Range("E13").Value = '  +9.01%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.466'
$ws.Range("E14").Value = '  +0.82%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '1.788.33'
$ws.Range("E15").Value = '  +3.57%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '7.021'
$ws.Range("E16").Value = '  -0.62%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001093'
$ws.Range("E17").Value = '  +1.55%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06647'
$ws.Range("E18").Value = '  -0.09%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '84.41'
$ws.Range("E19").Value = '  +2.86%  '

$ws.Range("E20").Value = '  +0.14%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.34'
$ws.Range("E21").Value = '  +4.68%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.473'
$ws.Range("E22").Value = '  +4.93%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '27.358.08'
$ws.Range("E23").Value = '  +3.16%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.53'
$ws.Range("E24").Value = '  -1.94%  '

$ws.Range("E25").Value = '  -0.16%  '

$ws.Range("B26").Value = 'ImmutableX'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.506'
$ws.Range("E26").Value = '  +5.97%  '

$ws.Range("B27").Value = 'LidoDAOToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.556'
$ws.Range("E27").Value = '  +6.37%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '21.32'
$ws.Range("E28").Value = '  +9.59%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '150.56'
$ws.Range("E29").Value = '  -0.41%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.989.71'
$ws.Range("E30").Value = '  +3.74%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '133.03'
$ws.Range("E31").Value = '  +1.20%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.061'
$ws.Range("E32").Value = '  -0.95%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.118'
$ws.Range("E33").Value = '  +1.99%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08689'
$ws.Range("E34").Value = '  +0.40%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '13.23'
$ws.Range("E35").Value = '  +3.37%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.673'
$ws.Range("E36").Value = '  -1.47%  '

$ws.Range("B37").Value = 'InternetComputer(DFINITY)'
$ws.Range("C37").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.424'
$ws.Range("E37").Value = '  +0.89%  '

$ws.Range("B38").Value = 'TheSandbox'
$ws.Range("C38").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.6874'
$ws.Range("E38").Value = '  +10.59%  '

$ws.Range("B39").Value = 'FraxShare'
$ws.Range("C39").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.836'
$ws.Range("E39").Value = '  +4.87%  '

$ws.Range("B40").Value = 'Hedera'
$ws.Range("C40").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.06337'
$ws.Range("E40").Value = '  +2.01%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.2192'
$ws.Range("E41").Value = '  +1.66%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.02339'
$ws.Range("E42").Value = '  +0.05%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.271'

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '14.45'
$ws.Range("E44").Value = '  +1.24%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.001'
$ws.Range("E45").Value = '  +0.09%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.6427'
$ws.Range("E46").Value = '  +6.44%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.848'
$ws.Range("E47").Value = '  +0.11%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.112'
$ws.Range("E48").Value = '  +3.04%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '129.36'
$ws.Range("E49").Value = '  +0.29%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07187'
$ws.Range("E50").Value = '  +0.07%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '78.99'
$ws.Range("E51").Value = '  +2.46%  '
